# NotifyTestData.xlsx fix: correct the QUERYSTRING values that were left
# over from copy/pasted test rows (OPQA-1013 "aggregated appreciation"
# tests and friends) so they use the same paging size ("?size=10") as the
# rest of the search-based dependency tests instead of the stray
# "?size=1" / "?size=2" / "?size=3" placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G16").Value = "?size=10"
$ws.Range("G22").Value = "?size=10"
$ws.Range("G23").Value = "?size=10"
$ws.Range("G31").Value = "?size=10"
$ws.Range("G32").Value = "?size=10"
$ws.Range("G33").Value = "?size=10"
$ws.Range("G39").Value = "?size=10"

# Restore the saved view/selection state recorded for the sheet.
$ws.Range("G1").Select()
$excel.ActiveWindow.ScrollColumn = 6
